$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels: B1/C1 were "X"/"Y", now become "Edad"/"Sexo",
# and D1 ("DIMENSION") value is preserved/moved to the end of the shared-string list.
$ws.Range("B1").Value = "Edad"
$ws.Range("C1").Value = "Sexo"
$ws.Range("D1").Value = "DIMENSION"

# Move the active selection to C1 (was D2).
$ws.Range("C1").Select()
